$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 1088:1089, pushing the existing data (previously at
# rows 1088-1208) down to rows 1090-1210. This matches the diff, which shows
# every row from 1090 onward now holding what used to be two rows above it,
# and the last two rows (1209, 1210) being new duplicates of the former
# 1207/1208 tail rows.
$ws.Rows("1088:1089").Insert()

# Row 1088 (new record) - a "Primera" quality entry for the same
# market/region/product, dated 45194. Rows().Insert() only pulls formatting
# down from the row above, not values, so every column must be written.
$ws.Cells.Item(1088, 1).Value2 = 6                                              # A - Mercado ID
$ws.Cells.Item(1088, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"    # B - Mercado
$ws.Cells.Item(1088, 3).Value2 = "Metropolitana"                               # C - Region
$ws.Cells.Item(1088, 4).Value2 = 45194                                          # D - Fecha
$ws.Cells.Item(1088, 5).Value2 = 13                                             # E - Codreg
$ws.Cells.Item(1088, 6).Value2 = 100112017                                      # F - Categoria ID
$ws.Cells.Item(1088, 7).Value2 = "Apio"                                         # G - Categoria
$ws.Cells.Item(1088, 8).Value2 = "Americana (o)"                                # H - Variedad
$ws.Cells.Item(1088, 9).Value2 = "Primera"                                      # I - Calidad
$ws.Cells.Item(1088, 10).Value2 = 1000                                          # J - Volumen
$ws.Cells.Item(1088, 11).Value2 = 6000                                          # K - Precio minimo
$ws.Cells.Item(1088, 12).Value2 = 6000                                          # L - Precio maximo
$ws.Cells.Item(1088, 13).Value2 = 6000                                          # M - Precio promedio ponderado
$ws.Cells.Item(1088, 14).Value2 = "`$/docena de matas"                          # N - Unidad de comercializacion
$ws.Cells.Item(1088, 15).Value2 = "Región de Coquimbo"                          # O - Origen
$ws.Cells.Item(1088, 16).Value2 = 1000                                          # P - Precio $/Kg
$ws.Cells.Item(1088, 17).Value2 = 6                                             # Q - Kg o Unidades
$ws.Cells.Item(1088, 18).Value2 = "Hortaliza"                                   # R - Clasificacion

# Row 1089 (new record) - the matching "Segunda" quality entry.
$ws.Cells.Item(1089, 1).Value2 = 6
$ws.Cells.Item(1089, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1089, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1089, 4).Value2 = 45194
$ws.Cells.Item(1089, 5).Value2 = 13
$ws.Cells.Item(1089, 6).Value2 = 100112017
$ws.Cells.Item(1089, 7).Value2 = "Apio"
$ws.Cells.Item(1089, 8).Value2 = "Americana (o)"
$ws.Cells.Item(1089, 9).Value2 = "Segunda"
$ws.Cells.Item(1089, 10).Value2 = 2600
$ws.Cells.Item(1089, 11).Value2 = 4500
$ws.Cells.Item(1089, 12).Value2 = 5500
$ws.Cells.Item(1089, 13).Value2 = 5077
$ws.Cells.Item(1089, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(1089, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(1089, 16).Value2 = 846
$ws.Cells.Item(1089, 17).Value2 = 6
$ws.Cells.Item(1089, 18).Value2 = "Hortaliza"
